$p = $ppt.ActivePresentation

# "STARS poster and document" - remove the poster slide (slide 3), which
# only contained a single full-slide picture (the STARS poster image).
# Deleting it also drops its now-unused slide-id entry from the
# presentation's slide list (and the related package plumbing: the
# slide's Content_Types override and its .rels part) automatically.
$p.Slides.Item(3).Delete()
